$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 288, shifting existing rows 288-325 down to 289-326.
$ws.Rows.Item(288).Insert()

# Populate the newly inserted row 288 with the new weekly data point.
$ws.Cells.Item(288,1).Value = 10
$ws.Cells.Item(288,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(288,3).Value = "La Araucanía"
$ws.Cells.Item(288,4).Value = 45124
$ws.Cells.Item(288,5).Value = 9
$ws.Cells.Item(288,6).Value = 100114007
$ws.Cells.Item(288,7).Value = "Jengibre"
$ws.Cells.Item(288,8).Value = "Sin especificar"
$ws.Cells.Item(288,9).Value = "Primera"
$ws.Cells.Item(288,10).Value = 180
$ws.Cells.Item(288,11).Value = 22000
$ws.Cells.Item(288,12).Value = 24000
$ws.Cells.Item(288,13).Value = 23111
$ws.Cells.Item(288,14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(288,15).Value = "Perú"
$ws.Cells.Item(288,16).Value = 1778
$ws.Cells.Item(288,17).Value = 13
$ws.Cells.Item(288,18).Value = "Hortaliza"

# Ensure the date column keeps the same number-format style used by the rest of column D.
$ws.Cells.Item(288,4).NumberFormat = $ws.Cells.Item(289,4).NumberFormat
